# Fix OOXMLValidator "Sch_UnexpectedElementContentExpectingComplex" error
# on several KeywordTok-like character styles in styles.xml.
#
# wml.xsd (CT_RPr) expects <w:b/>/<w:i/> to come before <w:color/>, but the
# styles below had them emitted in the wrong order, e.g.:
#   <w:rPr><w:color w:val="007020"/><w:b/></w:rPr>
# instead of:
#   <w:rPr><w:b/><w:color w:val="007020"/></w:rPr>
#
# Re-applying the existing Bold/Italic values through the Style.Font object
# model makes Word rewrite <w:rPr> in the schema-correct child order.

$d = $word.ActiveDocument

# styleId -> whether Bold / Italic should be (re-)applied, matching the
# diff: KeywordTok/ImportTok/ControlFlowTok/AlertTok/ErrorTok only need
# Bold moved before Color; CommentTok/DocumentationTok only need Italic
# moved before Color; AnnotationTok/CommentVarTok/InformationTok/WarningTok
# need both Bold and Italic moved before Color.
$fixes = @(
    @{ Name = "KeywordTok";       Bold = $true;  Italic = $false },
    @{ Name = "ImportTok";        Bold = $true;  Italic = $false },
    @{ Name = "CommentTok";       Bold = $false; Italic = $true  },
    @{ Name = "DocumentationTok"; Bold = $false; Italic = $true  },
    @{ Name = "AnnotationTok";    Bold = $true;  Italic = $true  },
    @{ Name = "CommentVarTok";    Bold = $true;  Italic = $true  },
    @{ Name = "ControlFlowTok";   Bold = $true;  Italic = $false },
    @{ Name = "InformationTok";   Bold = $true;  Italic = $true  },
    @{ Name = "WarningTok";       Bold = $true;  Italic = $true  },
    @{ Name = "AlertTok";         Bold = $true;  Italic = $false },
    @{ Name = "ErrorTok";         Bold = $true;  Italic = $false }
)

foreach ($fix in $fixes) {
    $s = $d.Styles($fix.Name)
    $f = $s.Font
    if ($fix.Bold) {
        $f.Bold = $true
    }
    if ($fix.Italic) {
        $f.Italic = $true
    }
}
